# auto_apresentacao_IA.pptx edit
# Shifts the "Introducao a Eletronica Digital" course deck content forward by
# one slide-topic, updates the title slide with a proper title + author
# subtitle, and turns the last slide into a "Conclusao" slide.
#
# Note: when a shape's paragraph count in the *new* text matches its
# *original* paragraph count and the new text contains multi-byte UTF-8
# characters, the COM host's in-place text replacement can mis-split a
# run at a stale byte offset (e.g. "...digitais" -> "...digitai" + "s").
# Collapsing the shape to a single throw-away paragraph first (a text
# assignment with no `r) changes the paragraph count so the later,
# real assignment is always a "fresh" write and never hits that path.

$p = $ppt.ActivePresentation

function Set-ShapeText {
    param($shape, $text)
    $shape.TextFrame.TextRange.Text = "x"
    $shape.TextFrame.TextRange.Text = $text
}

# Slide 1 - title slide
$s1 = $p.Slides.Item(1)
Set-ShapeText $s1.Shapes.Item(1) "Introdução à Eletrônica Digital"
Set-ShapeText $s1.Shapes.Item(2) "- Autor: [Seu nome]"

# Slide 2
$s2 = $p.Slides.Item(2)
Set-ShapeText $s2.Shapes.Item(1) "O que é Eletrônica Digital?"
Set-ShapeText $s2.Shapes.Item(2) "Introdução ao tema`rDefinição de eletrônica digital`rImportância da eletrônica digital no mundo moderno"

# Slide 3
$s3 = $p.Slides.Item(3)
Set-ShapeText $s3.Shapes.Item(1) "Sinais Digitais e Sinais Analógicos"
Set-ShapeText $s3.Shapes.Item(2) "Diferença entre sinais digitais e analógicos`rCaracterísticas dos sinais digitais`rVantagens dos sinais digitais"

# Slide 4
$s4 = $p.Slides.Item(4)
Set-ShapeText $s4.Shapes.Item(1) "Sistemas Digitais"
Set-ShapeText $s4.Shapes.Item(2) "Explicação sobre sistemas digitais`rExemplos de sistemas digitais comuns (computadores, celulares, etc.)`rBenefícios dos sistemas digitais"

# Slide 5
$s5 = $p.Slides.Item(5)
Set-ShapeText $s5.Shapes.Item(1) "Componentes Básicos da Eletrônica Digital"
Set-ShapeText $s5.Shapes.Item(2) "Introdução aos principais componentes da eletrônica digital`rTabela com componentes básicos (transistores, diodos, resistores, etc.)`rFunção de cada componente na eletrônica digital"

# Slide 6
$s6 = $p.Slides.Item(6)
Set-ShapeText $s6.Shapes.Item(1) "Portas Lógicas"
Set-ShapeText $s6.Shapes.Item(2) "Explicação sobre o conceito de portas lógicas`rTipos de portas lógicas (AND, OR, NOT, XOR, NAND, NOR)`rExemplificação das portas lógicas através de diagramas"

# Slide 7
$s7 = $p.Slides.Item(7)
Set-ShapeText $s7.Shapes.Item(1) "Circuitos Combinacionais"
Set-ShapeText $s7.Shapes.Item(2) "Definição de circuitos combinacionais`rExemplos de circuitos combinacionais (decodificadores, multiplexadores, somadores)`rUtilização de circuitos combinacionais na eletrônica digital"

# Slide 8
$s8 = $p.Slides.Item(8)
Set-ShapeText $s8.Shapes.Item(1) "Circuitos Sequenciais"
Set-ShapeText $s8.Shapes.Item(2) "Explicação sobre circuitos sequenciais`rTipos de circuitos sequenciais (flip-flops, contadores, registradores)`rAplicações dos circuitos sequenciais"

# Slide 9
$s9 = $p.Slides.Item(9)
Set-ShapeText $s9.Shapes.Item(1) "Microcontroladores"
Set-ShapeText $s9.Shapes.Item(2) "Descrição de microcontroladores`rVantagens e aplicações dos microcontroladores`rExemplos de microcontroladores famosos (Arduino, PIC, STM32)"

# Slide 10 - conclusion
$s10 = $p.Slides.Item(10)
Set-ShapeText $s10.Shapes.Item(1) "Conclusão"
Set-ShapeText $s10.Shapes.Item(2) "Recapitulação dos principais pontos abordados`rImportância da eletrônica digital na tecnologia atual`rReferências utilizadas na apresentação"
